# Apply updated 2018 population projection values to column B (rows 2-35)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prognosis")

$ws.Range("B2").Value = 11398589
$ws.Range("B3").Value = 7050034
$ws.Range("B4").Value = 10610055
$ws.Range("B5").Value = 5781190
$ws.Range("B6").Value = 82792351
$ws.Range("B7").Value = 9762334.200000018
$ws.Range("B8").Value = 10741165
$ws.Range("B9").Value = 46658447
$ws.Range("B10").Value = 122103780.1930661
$ws.Range("B11").Value = 4105493
$ws.Range("B12").Value = 60483973
$ws.Range("B13").Value = 1934379
$ws.Range("B14").Value = 602005
$ws.Range("B15").Value = 12712971.1031746
$ws.Range("B16").Value = 17181084
$ws.Range("B17").Value = 8822267
$ws.Range("B18").Value = 52375422.33903134
$ws.Range("B19").Value = 10291027
$ws.Range("B20").Value = 19530631
$ws.Range("B21").Value = 2066880
$ws.Range("B22").Value = 5443120
$ws.Range("B23").Value = 5513130
$ws.Range("B24").Value = 10120242
$ws.Range("B25").Value = 70634115.31
$ws.Range("B26").Value = 5234532.634779379
$ws.Range("B27").Value = 8484130
$ws.Range("B28").Value = 622359
$ws.Range("B29").Value = 2075108
$ws.Range("B30").Value = 2870324
$ws.Range("B31").Value = 7001444
$ws.Range("B32").Value = 3843183
$ws.Range("B33").Value = 348450
$ws.Range("B34").Value = 5658982.928571418
$ws.Range("B35").Value = 1319133

$wb.Save()
